$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows appended to the Indiana hospital ventilator dataset:
# 2021-03-31 through 2021-04-12 (rows 401-413), columns A (DATE) through O.
$newRows = @(
    @(401, "2021-03-31", 2191,121,1374,696,2787,45,479,2263,5.52,62.71,31.77,1.61,17.19,81.2),
    @(402, "2021-04-01", 2182,125,1327,730,2785,43,484,2258,5.73,60.82,33.46,1.54,17.38,81.08),
    @(403, "2021-04-02", 2181,133,1307,741,2781,45,489,2247,6.1,59.93,33.98,1.62,17.58,80.8),
    @(404, "2021-04-03", 2174,124,1256,794,2780,41,497,2242,5.7,57.77,36.52,1.47,17.88,80.65000000000001),
    @(405, "2021-04-04", 2172,123,1233,816,2777,38,484,2255,5.66,56.77,37.57,1.37,17.43,81.2),
    @(406, "2021-04-05", 2179,118,1204,857,2775,49,497,2229,5.42,55.25,39.33,1.77,17.91,80.31999999999999),
    @(407, "2021-04-06", 2176,144,1271,761,2775,62,515,2198,6.62,58.41,34.97,2.23,18.56,79.20999999999999),
    @(408, "2021-04-07", 2188,144,1364,680,2778,61,525,2192,6.58,62.34,31.08,2.2,18.9,78.91),
    @(409, "2021-04-08", 2191,156,1343,692,2772,58,508,2206,7.12,61.3,31.58,2.09,18.33,79.58),
    @(410, "2021-04-09", 2192,160,1334,698,2781,63,498,2220,7.3,60.86,31.84,2.27,17.91,79.83),
    @(411, "2021-04-10", 2205,171,1335,699,2785,69,488,2228,7.76,60.54,31.7,2.48,17.52,80),
    @(412, "2021-04-11", 2269,176,1310,783,2784,74,485,2225,7.76,57.73,34.51,2.66,17.42,79.92),
    @(413, "2021-04-12", 2187,178,1269,740,2781,79,481,2221,8.140000000000001,58.02,33.84,2.84,17.3,79.86)
)

# Excel auto-detects YYYY-MM-DD-looking text as a date and would otherwise store
# column A as a date serial number instead of the literal string. Pre-formatting
# the destination cells as Text ("@") forces the DATE column to stay a plain string,
# matching the existing rows (2..400) which are all shared-string text.
$dateRange = $ws.Range("A401:A413")
$dateRange.NumberFormat = "@"

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")
foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $row[$i + 2]
    }
}

# Drop the temporary Text formatting so the new DATE cells match the unstyled look
# of the existing rows (no explicit style index), then let the workbook recalc.
$dateRange.ClearFormats()
